$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 22
for ($block = 0; $block -lt 4; $block++) {
    for ($v = 1; $v -le 5; $v++) {
        $ws.Cells.Item($row, 1).Value = $v
        $ws.Cells.Item($row, 2).Value = $v
        $row++
    }
}
